$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.75736873280669
$ws.Range("D2").Value = 11.0651813289378
$ws.Range("E2").Value = 17.9515500816985
$ws.Range("F2").Value = 34.00903894237062
$ws.Range("G2").Value = 34.29816841915624
$ws.Range("H2").Value = 15.81602637450376
$ws.Range("J2").Value = 12.62912984307422
$ws.Range("K2").Value = 9.767162855987181
$ws.Range("L2").Value = 7.924163057567034
$ws.Range("M2").Value = 14.59644489691356
$ws.Range("O2").Value = 24.7841024839572

$ws.Range("B3").Value = 15.65967858411268
$ws.Range("D3").Value = 11.08111636696271
$ws.Range("E3").Value = 18.00430406571957
$ws.Range("F3").Value = 34.11796858626307
$ws.Range("G3").Value = 34.41861897301779
$ws.Range("H3").Value = 15.86577133542626
$ws.Range("J3").Value = 12.6559377599864
$ws.Range("K3").Value = 9.498155289324901
$ws.Range("L3").Value = 7.884784927224814
$ws.Range("M3").Value = 14.55928190945209
$ws.Range("O3").Value = 24.87103771344442

$ws.Range("B4").Value = 15.60182502278972
$ws.Range("D4").Value = 11.09228595307725
$ws.Range("E4").Value = 18.03857771411939
$ws.Range("F4").Value = 34.19145638514679
$ws.Range("G4").Value = 34.50135875395598
$ws.Range("H4").Value = 15.89848928643692
$ws.Range("J4").Value = 12.67326564719493
$ws.Range("K4").Value = 9.327735112853187
$ws.Range("L4").Value = 7.861074376584336
$ws.Range("M4").Value = 14.53793255374607
$ws.Range("O4").Value = 24.92886547676805

$ws.Range("B5").Value = 15.5788037990011
$ws.Range("D5").Value = 11.0971865291238
$ws.Range("E5").Value = 18.05301899352087
$ws.Range("F5").Value = 34.22306251027188
$ws.Range("G5").Value = 34.53727839839515
$ws.Range("H5").Value = 15.91236935422583
$ws.Range("J5").Value = 12.68054573148009
$ws.Range("K5").Value = 9.257035510853274
$ws.Range("L5").Value = 7.851535487424092
$ws.Range("M5").Value = 14.52960725816567
$ws.Range("O5").Value = 24.95354902496967

$ws.Range("B6").Value = 15.5750151860152
$ws.Range("D6").Value = 11.09802135168146
$ws.Range("E6").Value = 18.05544564803852
$ws.Range("F6").Value = 34.22841085177517
$ws.Range("G6").Value = 34.543375674522
$ws.Range("H6").Value = 15.91470719493364
$ws.Range("J6").Value = 12.6817678193862
$ws.Range("K6").Value = 9.245222363377424
$ws.Range("L6").Value = 7.849959176690859
$ws.Range("M6").Value = 14.52824763477997
$ws.Range("O6").Value = 24.95771523450992

$ws.Range("B7").Value = 15.60151228014051
$ws.Range("D7").Value = 11.09235063070642
$ws.Range("E7").Value = 18.03877055157217
$ws.Range("F7").Value = 34.19187592023154
$ws.Range("G7").Value = 34.50183426904427
$ws.Range("H7").Value = 15.89867426150136
$ws.Range("J7").Value = 12.67336294207878
$ws.Range("K7").Value = 9.326786609131869
$ws.Range("L7").Value = 7.860945224886677
$ws.Range("M7").Value = 14.53781875166974
$ws.Range("O7").Value = 24.92919384061033

$ws.Range("B8").Value = 15.72325520872539
$ws.Range("D8").Value = 11.07038850147232
$ws.Range("E8").Value = 17.96934962102891
$ws.Range("F8").Value = 34.04522622306725
$ws.Range("G8").Value = 34.33787343505153
$ws.Range("H8").Value = 15.8327274414082
$ws.Range("J8").Value = 12.63819350095969
$ws.Range("K8").Value = 9.675537362302993
$ws.Range("L8").Value = 7.910491866483839
$ws.Range("M8").Value = 14.58332953478005
$ws.Range("O8").Value = 24.8131540188447

$ws.Range("B9").Value = 15.9779989736982
$ws.Range("D9").Value = 11.03829043112509
$ws.Range("E9").Value = 17.84810099288238
$ws.Range("F9").Value = 33.81011573069668
$ws.Range("G9").Value = 34.08628776413008
$ws.Range("H9").Value = 15.72063668482282
$ws.Range("J9").Value = 12.57608221249171
$ws.Range("K9").Value = 10.31508756760121
$ws.Range("L9").Value = 8.011107745867193
$ws.Range("M9").Value = 14.68398551460932
$ws.Range("O9").Value = 24.62092506846079

$ws.Range("B10").Value = 16.17369049969123
$ws.Range("D10").Value = 11.02136330775364
$ws.Range("E10").Value = 17.76802293126565
$ws.Range("F10").Value = 33.66945075757997
$ws.Range("G10").Value = 33.94441576628179
$ws.Range("H10").Value = 15.64875793837363
$ws.Range("J10").Value = 12.53458773393591
$ws.Range("K10").Value = 10.75467961411015
$ws.Range("L10").Value = 8.086786341916351
$ws.Range("M10").Value = 14.76455152388369
$ws.Range("O10").Value = 24.50125925805748

$ws.Range("B11").Value = 16.2642958487151
$ws.Range("D11").Value = 11.01510035727279
$ws.Range("E11").Value = 17.73353323418264
$ws.Range("F11").Value = 33.6124384675831
$ws.Range("G11").Value = 33.88926825894135
$ws.Range("H11").Value = 15.61832662093217
$ws.Range("J11").Value = 12.51660081398096
$ws.Range("K11").Value = 10.94749959297156
$ws.Range("L11").Value = 8.121519136416154
$ws.Range("M11").Value = 14.80256178460855
$ws.Range("O11").Value = 24.4515089975286

$ws.Range("B12").Value = 16.29881001512496
$ws.Range("D12").Value = 11.01293472106687
$ws.Range("E12").Value = 17.72075046745429
$ws.Range("F12").Value = 33.5918538376056
$ws.Range("G12").Value = 33.86974034089504
$ws.Range("H12").Value = 15.60712858560592
$ws.Range("J12").Value = 12.50991683884248
$ws.Range("K12").Value = 11.0194443674152
$ws.Range("L12").Value = 8.134708942576548
$ws.Range("M12").Value = 14.81714394662271
$ws.Range("O12").Value = 24.43334431569587

$ws.Range("B13").Value = 16.29136808764272
$ws.Range("D13").Value = 11.0133919784948
$ws.Range("E13").Value = 17.72349112821554
$ws.Range("F13").Value = 33.59624240115988
$ws.Range("G13").Value = 33.87388567959365
$ws.Range("H13").Value = 15.60952580449361
$ws.Range("J13").Value = 12.51135070070862
$ws.Range("K13").Value = 11.00399809300501
$ws.Range("L13").Value = 8.13186673368334
$ws.Range("M13").Value = 14.81399514900975
$ws.Range("O13").Value = 24.43722638705392

$ws.Range("B14").Value = 16.26713140173858
$ws.Range("D14").Value = 11.01491806438492
$ws.Range("E14").Value = 17.73247602765882
$ws.Range("F14").Value = 33.6107248146351
$ws.Range("G14").Value = 33.8876344978696
$ws.Range("H14").Value = 15.617398827137
$ws.Range("J14").Value = 12.51604837169173
$ws.Range("K14").Value = 10.95344024952554
$ws.Range("L14").Value = 8.122603562199625
$ws.Range("M14").Value = 14.80375773151305
$ws.Range("O14").Value = 24.45000105262715

$ws.Range("B15").Value = 16.25231156924141
$ws.Range("D15").Value = 11.0158796436206
$ws.Range("E15").Value = 17.73801567815588
$ws.Range("F15").Value = 33.61972659066709
$ws.Range("G15").Value = 33.89623267622289
$ws.Range("H15").Value = 15.62226368047449
$ws.Range("J15").Value = 12.51894238961308
$ws.Range("K15").Value = 10.9223312965099
$ws.Range("L15").Value = 8.116934251051262
$ws.Range("M15").Value = 14.79751135132831
$ws.Range("O15").Value = 24.45791378770311

$ws.Range("B16").Value = 16.16779906751431
$ws.Range("D16").Value = 11.02180147157463
$ws.Range("E16").Value = 17.77031583189281
$ws.Range("F16").Value = 33.67331716216798
$ws.Range("G16").Value = 33.94820917740747
$ws.Range("H16").Value = 15.65079228200231
$ws.Range("J16").Value = 12.53578106374741
$ws.Range("K16").Value = 10.74193063036141
$ws.Range("L16").Value = 8.084522055537297
$ws.Range("M16").Value = 14.76209423437612
$ws.Range("O16").Value = 24.50460492205586

$ws.Range("B17").Value = 16.11634222324991
$ws.Range("D17").Value = 11.02580197855183
$ws.Range("E17").Value = 17.79062665243848
$ws.Range("F17").Value = 33.70798124580076
$ws.Range("G17").Value = 33.98250387410094
$ws.Range("H17").Value = 15.6688739897111
$ws.Range("J17").Value = 12.54633836129363
$ws.Range("K17").Value = 10.62939636194992
$ws.Range("L17").Value = 8.064711847015246
$ws.Range("M17").Value = 14.74071040665106
$ws.Range("O17").Value = 24.53444925753239

$ws.Range("B18").Value = 16.08689618128267
$ws.Range("D18").Value = 11.02823829104755
$ws.Range("E18").Value = 17.80249139103923
$ws.Range("F18").Value = 33.72857573344555
$ws.Range("G18").Value = 34.00311291599417
$ws.Range("H18").Value = 15.67948747840849
$ws.Range("J18").Value = 12.55249436240067
$ws.Range("K18").Value = 10.56399868589796
$ws.Range("L18").Value = 8.053346772266432
$ws.Range("M18").Value = 14.72853947941154
$ws.Range("O18").Value = 24.5520559832113

$ws.Range("B19").Value = 16.07695284325326
$ws.Range("D19").Value = 11.02908644689028
$ws.Range("E19").Value = 17.80653996085468
$ws.Range("F19").Value = 33.73566140352928
$ws.Range("G19").Value = 34.01024239432699
$ws.Range("H19").Value = 15.68311767893962
$ws.Range("J19").Value = 12.55459307565414
$ws.Range("K19").Value = 10.5417423189797
$ws.Range("L19").Value = 8.049503979778748
$ws.Range("M19").Value = 14.72444089350285
$ws.Range("O19").Value = 24.55809304534668

$ws.Range("B20").Value = 16.12180448083828
$ws.Range("D20").Value = 11.02536211669591
$ws.Range("E20").Value = 17.78844565083837
$ws.Range("F20").Value = 33.7042232261666
$ws.Range("G20").Value = 33.97876165515064
$ws.Range("H20").Value = 15.66692708041336
$ws.Range("J20").Value = 12.54520585812924
$ws.Range("K20").Value = 10.64144559565579
$ws.Range("L20").Value = 8.066817700955644
$ws.Range("M20").Value = 14.7429735048328
$ws.Range("O20").Value = 24.53122662402179

$ws.Range("B21").Value = 16.2742449624416
$ws.Range("D21").Value = 11.01446423055262
$ws.Range("E21").Value = 17.72982941505288
$ws.Range("F21").Value = 33.60644370121274
$ws.Range("G21").Value = 33.88355931901894
$ws.Range("H21").Value = 15.61507749402444
$ws.Range("J21").Value = 12.51466510202906
$ws.Range("K21").Value = 10.96831973048388
$ws.Range("L21").Value = 8.125323423569661
$ws.Range("M21").Value = 14.80675965107175
$ws.Range("O21").Value = 24.44623050748266

$ws.Range("B22").Value = 16.37505078066674
$ws.Range("D22").Value = 11.00854225843151
$ws.Range("E22").Value = 17.69313871058567
$ws.Range("F22").Value = 33.54839577357478
$ws.Range("G22").Value = 33.8292404528278
$ws.Range("H22").Value = 15.58308876046764
$ws.Range("J22").Value = 12.4954465812332
$ws.Range("K22").Value = 11.17568558363343
$ws.Range("L22").Value = 8.163774239158213
$ws.Range("M22").Value = 14.8495428711194
$ws.Range("O22").Value = 24.39461342228363

$ws.Range("B23").Value = 16.32114907606377
$ws.Range("D23").Value = 11.01159332012757
$ws.Range("E23").Value = 17.71257346767609
$ws.Range("F23").Value = 33.57884072916457
$ws.Range("G23").Value = 33.85750698246243
$ws.Range("H23").Value = 15.59998818324134
$ws.Range("J23").Value = 12.50563619674503
$ws.Range("K23").Value = 11.06559700358247
$ws.Range("L23").Value = 8.143234993812511
$ws.Range("M23").Value = 14.82661084663844
$ws.Range("O23").Value = 24.42180233125853

$ws.Range("B24").Value = 16.11933456505633
$ws.Range("D24").Value = 11.02556055343189
$ws.Range("E24").Value = 17.78943109668704
$ws.Range("F24").Value = 33.70592015352779
$ws.Range("G24").Value = 33.98045073244555
$ws.Range("H24").Value = 15.6678065987959
$ws.Range("J24").Value = 12.54571759342033
$ws.Range("K24").Value = 10.63600031482731
$ws.Range("L24").Value = 8.065865568876339
$ws.Range("M24").Value = 14.74194997472627
$ws.Range("O24").Value = 24.53268217864763

$ws.Range("B25").Value = 15.9074947489101
$ws.Range("D25").Value = 11.04580259848189
$ws.Range("E25").Value = 17.8793157715756
$ws.Range("F25").Value = 33.86809305019259
$ws.Range("G25").Value = 34.14682523385673
$ws.Range("H25").Value = 15.74911866719374
$ws.Range("J25").Value = 12.59215523403276
$ws.Range("K25").Value = 10.14717560363419
$ws.Range("L25").Value = 7.983554284177336
$ws.Range("M25").Value = 14.65556821992005
$ws.Range("O25").Value = 24.66914319157414
